$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the loss-function labels in column A.
# Assignment order matters for how the shared-strings table gets rebuilt,
# so we set them in the order that preserves the original string slots:
#   slot 15 "Ours Loss"             -> "Our Loss"              (A6)
#   slot 16 "Image-Text-Contrastive" -> "Image-Text Contrastive" (A4)
#   slot 17 "UniC Loss"             -> "UniCL Loss"            (A5)
$ws.Range("A6").Value = "Our Loss"
$ws.Range("A4").Value = "Image-Text Contrastive"
$ws.Range("A5").Value = "UniCL Loss"
